$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-28 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-29 Saturday", 2) | Out-Null
$d.Content.Find.Execute("70+10=80", $true, $false, $false, $false, $false, $true, 1, $false, "3+25=28", 2) | Out-Null
$d.Content.Find.Execute("74-32=42", $true, $false, $false, $false, $false, $true, 1, $false, "44-38=6", 2) | Out-Null
$d.Content.Find.Execute("73-21=52", $true, $false, $false, $false, $false, $true, 1, $false, "69+16=85", 2) | Out-Null
$d.Content.Find.Execute("72-34=38", $true, $false, $false, $false, $false, $true, 1, $false, "40-14=26", 2) | Out-Null
$d.Content.Find.Execute("93-45=48", $true, $false, $false, $false, $false, $true, 1, $false, "26-2=24", 2) | Out-Null
$d.Content.Find.Execute("42-26=16", $true, $false, $false, $false, $false, $true, 1, $false, "19+2=21", 2) | Out-Null
$d.Content.Find.Execute("70-41=29", $true, $false, $false, $false, $false, $true, 1, $false, "55-27=28", 2) | Out-Null
$d.Content.Find.Execute("13+27=40", $true, $false, $false, $false, $false, $true, 1, $false, "93-9=84", 2) | Out-Null
$d.Content.Find.Execute("69-43=26", $true, $false, $false, $false, $false, $true, 1, $false, "37+25=62", 2) | Out-Null
$d.Content.Find.Execute("98-75=23", $true, $false, $false, $false, $false, $true, 1, $false, "92-88=4", 2) | Out-Null
$d.Content.Find.Execute("12+32=44", $true, $false, $false, $false, $false, $true, 1, $false, "46+0=46", 2) | Out-Null
$d.Content.Find.Execute("35+49=84", $true, $false, $false, $false, $false, $true, 1, $false, "60-22=38", 2) | Out-Null
$d.Content.Find.Execute("60+7=67", $true, $false, $false, $false, $false, $true, 1, $false, "91-27=64", 2) | Out-Null
$d.Content.Find.Execute("72-59=13", $true, $false, $false, $false, $false, $true, 1, $false, "54-21=33", 2) | Out-Null
$d.Content.Find.Execute("89-82=7", $true, $false, $false, $false, $false, $true, 1, $false, "90-63=27", 2) | Out-Null
$d.Content.Find.Execute("73+4=77", $true, $false, $false, $false, $false, $true, 1, $false, "34+30=64", 2) | Out-Null
$d.Content.Find.Execute("99-19=80", $true, $false, $false, $false, $false, $true, 1, $false, "24+54=78", 2) | Out-Null
$d.Content.Find.Execute("45+52=97", $true, $false, $false, $false, $false, $true, 1, $false, "32+30=62", 2) | Out-Null
$d.Content.Find.Execute("97-71=26", $true, $false, $false, $false, $false, $true, 1, $false, "70+5=75", 2) | Out-Null
$d.Content.Find.Execute("46-14=32", $true, $false, $false, $false, $false, $true, 1, $false, "45-15=30", 2) | Out-Null
$d.Content.Find.Execute("61-35=26", $true, $false, $false, $false, $false, $true, 1, $false, "57-47=10", 2) | Out-Null
$d.Content.Find.Execute("45+8=53", $true, $false, $false, $false, $false, $true, 1, $false, "96-38=58", 2) | Out-Null
$d.Content.Find.Execute("66-41=25", $true, $false, $false, $false, $false, $true, 1, $false, "28+57=85", 2) | Out-Null
$d.Content.Find.Execute("65-62=3", $true, $false, $false, $false, $false, $true, 1, $false, "28+34=62", 2) | Out-Null
$d.Content.Find.Execute("26+13=39", $true, $false, $false, $false, $false, $true, 1, $false, "83-14=69", 2) | Out-Null
$d.Content.Find.Execute("12+19=31", $true, $false, $false, $false, $false, $true, 1, $false, "53+42=95", 2) | Out-Null
$d.Content.Find.Execute("49-18=31", $true, $false, $false, $false, $false, $true, 1, $false, "74-58=16", 2) | Out-Null
$d.Content.Find.Execute("99-27=72", $true, $false, $false, $false, $false, $true, 1, $false, "3+38=41", 2) | Out-Null
$d.Content.Find.Execute("54+34=88", $true, $false, $false, $false, $false, $true, 1, $false, "52-40=12", 2) | Out-Null
$d.Content.Find.Execute("1+10=11", $true, $false, $false, $false, $false, $true, 1, $false, "28-24=4", 2) | Out-Null
$d.Content.Find.Execute("7+69=76", $true, $false, $false, $false, $false, $true, 1, $false, "68+29=97", 2) | Out-Null
$d.Content.Find.Execute("30-30=0", $true, $false, $false, $false, $false, $true, 1, $false, "63+16=79", 2) | Out-Null
$d.Content.Find.Execute("22+70=92", $true, $false, $false, $false, $false, $true, 1, $false, "57-43=14", 2) | Out-Null
$d.Content.Find.Execute("22+40=62", $true, $false, $false, $false, $false, $true, 1, $false, "70-13=57", 2) | Out-Null
$d.Content.Find.Execute("34-5=29", $true, $false, $false, $false, $false, $true, 1, $false, "77-47=30", 2) | Out-Null
$d.Content.Find.Execute("42-15=27", $true, $false, $false, $false, $false, $true, 1, $false, "14+71=85", 2) | Out-Null
$d.Content.Find.Execute("6+80=86", $true, $false, $false, $false, $false, $true, 1, $false, "90-71=19", 2) | Out-Null
$d.Content.Find.Execute("11+25=36", $true, $false, $false, $false, $false, $true, 1, $false, "28-15=13", 2) | Out-Null
$d.Content.Find.Execute("33-20=13", $true, $false, $false, $false, $false, $true, 1, $false, "10+86=96", 2) | Out-Null
$d.Content.Find.Execute("27-24=3", $true, $false, $false, $false, $false, $true, 1, $false, "37+13=50", 2) | Out-Null
$d.Content.Find.Execute("54-33=21", $true, $false, $false, $false, $false, $true, 1, $false, "6+4=10", 2) | Out-Null
$d.Content.Find.Execute("79+11=90", $true, $false, $false, $false, $false, $true, 1, $false, "33+35=68", 2) | Out-Null
$d.Content.Find.Execute("89-79=10", $true, $false, $false, $false, $false, $true, 1, $false, "35+27=62", 2) | Out-Null
$d.Content.Find.Execute("54-40=14", $true, $false, $false, $false, $false, $true, 1, $false, "16+7=23", 2) | Out-Null
$d.Content.Find.Execute("80-69=11", $true, $false, $false, $false, $false, $true, 1, $false, "27-20=7", 2) | Out-Null
$d.Content.Find.Execute("77-51=26", $true, $false, $false, $false, $false, $true, 1, $false, "38+20=58", 2) | Out-Null
$d.Content.Find.Execute("6+33=39", $true, $false, $false, $false, $false, $true, 1, $false, "97-67=30", 2) | Out-Null
$d.Content.Find.Execute("39+24=63", $true, $false, $false, $false, $false, $true, 1, $false, "19+10=29", 2) | Out-Null
$d.Content.Find.Execute("35+63=98", $true, $false, $false, $false, $false, $true, 1, $false, "61+2=63", 2) | Out-Null
$d.Content.Find.Execute("52-30=22", $true, $false, $false, $false, $false, $true, 1, $false, "31+45=76", 2) | Out-Null
$d.Content.Find.Execute("1+90=91", $true, $false, $false, $false, $false, $true, 1, $false, "6-2=4", 2) | Out-Null
$d.Content.Find.Execute("97-9=88", $true, $false, $false, $false, $false, $true, 1, $false, "35-13=22", 2) | Out-Null
$d.Content.Find.Execute("88-18=70", $true, $false, $false, $false, $false, $true, 1, $false, "33+31=64", 2) | Out-Null
$d.Content.Find.Execute("65-22=43", $true, $false, $false, $false, $false, $true, 1, $false, "52+6=58", 2) | Out-Null
$d.Content.Find.Execute("48+10=58", $true, $false, $false, $false, $false, $true, 1, $false, "27-6=21", 2) | Out-Null
$d.Content.Find.Execute("65-43=22", $true, $false, $false, $false, $false, $true, 1, $false, "96-19=77", 2) | Out-Null
$d.Content.Find.Execute("16-8=8", $true, $false, $false, $false, $false, $true, 1, $false, "83-12=71", 2) | Out-Null
$d.Content.Find.Execute("56+21=77", $true, $false, $false, $false, $false, $true, 1, $false, "85-33=52", 2) | Out-Null
$d.Content.Find.Execute("9+45=54", $true, $false, $false, $false, $false, $true, 1, $false, "86-11=75", 2) | Out-Null
$d.Content.Find.Execute("91-86=5", $true, $false, $false, $false, $false, $true, 1, $false, "38+45=83", 2) | Out-Null
$d.Content.Find.Execute("74+13=87", $true, $false, $false, $false, $false, $true, 1, $false, "59+34=93", 2) | Out-Null
$d.Content.Find.Execute("70-68=2", $true, $false, $false, $false, $false, $true, 1, $false, "57-6=51", 2) | Out-Null
$d.Content.Find.Execute("85-31=54", $true, $false, $false, $false, $false, $true, 1, $false, "39-28=11", 2) | Out-Null
$d.Content.Find.Execute("66+13=79", $true, $false, $false, $false, $false, $true, 1, $false, "48-45=3", 2) | Out-Null
$d.Content.Find.Execute("25+22=47", $true, $false, $false, $false, $false, $true, 1, $false, "77-21=56", 2) | Out-Null
$d.Content.Find.Execute("24+25=49", $true, $false, $false, $false, $false, $true, 1, $false, "90+0=90", 2) | Out-Null
$d.Content.Find.Execute("53-14=39", $true, $false, $false, $false, $false, $true, 1, $false, "31+45=76", 2) | Out-Null
$d.Content.Find.Execute("63+19=82", $true, $false, $false, $false, $false, $true, 1, $false, "31+22=53", 2) | Out-Null
$d.Content.Find.Execute("38+4=42", $true, $false, $false, $false, $false, $true, 1, $false, "26-16=10", 2) | Out-Null
$d.Content.Find.Execute("22+22=44", $true, $false, $false, $false, $false, $true, 1, $false, "25+65=90", 2) | Out-Null
$d.Content.Find.Execute("64-56=8", $true, $false, $false, $false, $false, $true, 1, $false, "78+2=80", 2) | Out-Null
$d.Content.Find.Execute("12+84=96", $true, $false, $false, $false, $false, $true, 1, $false, "10+20=30", 2) | Out-Null
$d.Content.Find.Execute("43+32=75", $true, $false, $false, $false, $false, $true, 1, $false, "39+32=71", 2) | Out-Null
$d.Content.Find.Execute("25+46=71", $true, $false, $false, $false, $false, $true, 1, $false, "28+27=55", 2) | Out-Null
$d.Content.Find.Execute("14+45=59", $true, $false, $false, $false, $false, $true, 1, $false, "87-52=35", 2) | Out-Null
$d.Content.Find.Execute("20+32=52", $true, $false, $false, $false, $false, $true, 1, $false, "50+28=78", 2) | Out-Null
$d.Content.Find.Execute("76-69=7", $true, $false, $false, $false, $false, $true, 1, $false, "60+34=94", 2) | Out-Null
$d.Content.Find.Execute("18+73=91", $true, $false, $false, $false, $false, $true, 1, $false, "15+8=23", 2) | Out-Null
$d.Content.Find.Execute("27+37=64", $true, $false, $false, $false, $false, $true, 1, $false, "78-62=16", 2) | Out-Null
$d.Content.Find.Execute("10+39=49", $true, $false, $false, $false, $false, $true, 1, $false, "81-9=72", 2) | Out-Null
$d.Content.Find.Execute("85+12=97", $true, $false, $false, $false, $false, $true, 1, $false, "32+4=36", 2) | Out-Null
$d.Content.Find.Execute("94-29=65", $true, $false, $false, $false, $false, $true, 1, $false, "26-17=9", 2) | Out-Null
$d.Content.Find.Execute("1-0=1", $true, $false, $false, $false, $false, $true, 1, $false, "66-64=2", 2) | Out-Null
$d.Content.Find.Execute("92+2=94", $true, $false, $false, $false, $false, $true, 1, $false, "40-35=5", 2) | Out-Null
$d.Content.Find.Execute("46-8=38", $true, $false, $false, $false, $false, $true, 1, $false, "0+1=1", 2) | Out-Null
$d.Content.Find.Execute("72+23=95", $true, $false, $false, $false, $false, $true, 1, $false, "51-43=8", 2) | Out-Null
$d.Content.Find.Execute("83-5=78", $true, $false, $false, $false, $false, $true, 1, $false, "80-56=24", 2) | Out-Null
$d.Content.Find.Execute("7+30=37", $true, $false, $false, $false, $false, $true, 1, $false, "21-15=6", 2) | Out-Null
$d.Content.Find.Execute("74-37=37", $true, $false, $false, $false, $false, $true, 1, $false, "58-13=45", 2) | Out-Null
$d.Content.Find.Execute("9+81=90", $true, $false, $false, $false, $false, $true, 1, $false, "38-24=14", 2) | Out-Null
$d.Content.Find.Execute("50+25=75", $true, $false, $false, $false, $false, $true, 1, $false, "36+46=82", 2) | Out-Null
$d.Content.Find.Execute("9+66=75", $true, $false, $false, $false, $false, $true, 1, $false, "79-78=1", 2) | Out-Null
$d.Content.Find.Execute("30+1=31", $true, $false, $false, $false, $false, $true, 1, $false, "39+54=93", 2) | Out-Null
$d.Content.Find.Execute("99+0=99", $true, $false, $false, $false, $false, $true, 1, $false, "33+27=60", 2) | Out-Null
$d.Content.Find.Execute("57+42=99", $true, $false, $false, $false, $false, $true, 1, $false, "20-3=17", 2) | Out-Null
$d.Content.Find.Execute("90-50=40", $true, $false, $false, $false, $false, $true, 1, $false, "43-40=3", 2) | Out-Null
$d.Content.Find.Execute("51-26=25", $true, $false, $false, $false, $false, $true, 1, $false, "55+25=80", 2) | Out-Null
$d.Content.Find.Execute("49-33=16", $true, $false, $false, $false, $false, $true, 1, $false, "60+8=68", 2) | Out-Null
$d.Content.Find.Execute("79-31=48", $true, $false, $false, $false, $false, $true, 1, $false, "42+55=97", 2) | Out-Null
$d.Content.Find.Execute("40+41=81", $true, $false, $false, $false, $false, $true, 1, $false, "83-6=77", 2) | Out-Null
